{"js": "// Insert \" el dia xx de xxxx de xxxx\" right after the existing run of text\n// \"impartir tr\u00e1mite al memorial radicado\" (and before the following \", \"),\n// so the sentence reads \"...impartir tr\u00e1mite al memorial radicado el dia xx\n// de xxxx de xxxx, mediante el cual se present\u00f3 liquidaci\u00f3n de cr\u00e9dito.\"\n\nconst body = context.document.body;\nconst anchor = \"impartir tr\u00e1mite al memorial radicado\";\n\nconst results = body.search(anchor, { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(`Anchor text not found: ${anchor}`);\n}\n\n// There is a single occurrence of this phrase in the document.\nconst target = results.items[0];\n\n// Insert the new wording immediately after the matched phrase.\ntarget.insertText(\" el dia xx de xxxx de xxxx\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Insert \" el dia xx de xxxx de xxxx\" right after the existing text\n# \"impartir tr\u00e1mite al memorial radicado\" (and before the following \", \"),\n# so the sentence reads \"...impartir tr\u00e1mite al memorial radicado el dia xx\n# de xxxx de xxxx, mediante el cual se present\u00f3 liquidaci\u00f3n de cr\u00e9dito.\"\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"impartir tr\u00e1mite al memorial radicado\"\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWholeWord = $false\n$rng.Find.MatchWildcards = $false\n$rng.Find.Forward = $true\n\n$found = $rng.Find.Execute()\n\nif ($found) {\n    # $rng now spans exactly the matched phrase; InsertAfter drops the new\n    # text right after it, ahead of the following \", mediante el cual...\" run.\n    $rng.InsertAfter(\" el dia xx de xxxx de xxxx\")\n}\n\n$d.Save()\n"}
